$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# --- Row 4 ---
$ws.Range("C4").Value = "high"
$ws.Range("E4").Value = "Implementing game-renderer"
$ws.Range("G4").Value = 1
$ws.Range("H4").Value = "user"
$ws.Range("I4").Value = "neccessity for base game"

# --- Row 5 ---
$ws.Range("C5").Value = "high"
$ws.Range("E5").Value = "Implementing movable entities"
$ws.Range("G5").Value = 2
$ws.Range("H5").Value = "user"
$ws.Range("I5").Value = "neccessity for base game"

# --- Row 6 ---
$ws.Range("C6").Value = "high"
$ws.Range("E6").Value = "Implementing algorithm for pathfinding (A*)"
$ws.Range("G6").Value = 3
$ws.Range("H6").Value = "user"
$ws.Range("I6").Value = "I want the enemies to target my base and not have them move random"

# --- Row 7 ---
$ws.Range("C7").Value = "high"
$ws.Range("E7").Value = "Implementing stationary entities"
$ws.Range("G7").Value = 4
$ws.Range("H7").Value = "user"
$ws.Range("I7").Value = "neccessity for base game"

# --- Row 8 ---
$ws.Range("C8").Value = "high"
$ws.Range("E8").Value = "Implementing game class"
$ws.Range("G8").Value = 5
$ws.Range("H8").Value = "user"
$ws.Range("I8").Value = "neccessity for base game"

# --- Row 9 ---
$ws.Range("C9").Value = "medium"
$ws.Range("E9").Value = "Implementing database schema for save states"
$ws.Range("G9").Value = 6
$ws.Range("H9").Value = "user"
$ws.Range("I9").Value = "to save my level progress"

# --- Row 10 ---
$ws.Range("C10").Value = "medium"
$ws.Range("E10").Value = "Implementing saving gamestate in DB"
$ws.Range("G10").Value = 7
$ws.Range("H10").Value = "user"
$ws.Range("I10").Value = "to save my level progress"

# --- Row 11 ---
$ws.Range("C11").Value = "low"
$ws.Range("E11").Value = "Implementing 1 aditional tower"
$ws.Range("G11").Value = 8
$ws.Range("H11").Value = "user"
$ws.Range("I11").Value = "I want different towers to chose from"

# --- Row 12 ---
$ws.Range("C12").Value = "high"
$ws.Range("E12").Value = "Implementing Fog of War"
$ws.Range("G12").Value = 9
$ws.Range("H12").Value = "user"
$ws.Range("I12").Value = "I want to have to play to see the enemy base, and not from the beginning"

# --- Row 13 ---
$ws.Range("C13").Value = "high"
$ws.Range("E13").Value = "Implementing one static game field"
$ws.Range("G13").Value = 10
$ws.Range("H13").Value = "user"
$ws.Range("I13").Value = "neccessity for base game"

# --- Row 14 ---
$ws.Range("C14").Value = "low"
$ws.Range("E14").Value = "Implementing an algorithm for random game field generation"
$ws.Range("G14").Value = 11
$ws.Range("H14").Value = "user"
$ws.Range("I14").Value = "I want random generated game fields and not one static only"

# --- Row 15 ---
$ws.Range("C15").Value = "low"
$ws.Range("E15").Value = "save format for game field, for random generation as well as static fields, or maybe self edited"
$ws.Range("G15").Value = 12
$ws.Range("H15").Value = "user"
$ws.Range("I15").Value = $null

# --- Row 16 (cleared out) ---
$ws.Range("C16").Value = $null
$ws.Range("E16").Value = $null
$ws.Range("G16").Value = $null
$ws.Range("H16").Value = $null
$ws.Range("I16").Value = $null

# --- Wrap text / left-align formatting on the longer "Product Goal" cells ---
foreach ($r in 6,7,10,11,12,13,14,15,16) {
    $c = $ws.Range("E$r")
    $c.WrapText = $true
    $c.HorizontalAlignment = -4131
}

# --- Row heights for rows whose text now wraps onto multiple lines ---
$ws.Rows.Item(6).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(12).RowHeight = 45
$ws.Rows.Item(14).RowHeight = 45

# --- Data validation ranges now skip rows 15:16 ---
$ws.Range("C4:C37").Validation.Delete()
$ws.Range("C4:C14").Validation.Add(3, 1, 1, "=$D$1:$G$1")
$ws.Range("C17:C37").Validation.Add(3, 1, 1, "=$D$1:$G$1")

$ws.Range("D4:D37").Validation.Delete()
$ws.Range("D4:D14").Validation.Add(3, 1, 1, "=$D$2:$G$2")
$ws.Range("D17:D37").Validation.Add(3, 1, 1, "=$D$2:$G$2")

# --- Selection / view state ---
$ws.Range("I15").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
